# "Materials Done. Scenes and Meshes started."
# Fill in the next two days of logged time (rows 8 & 9) on the time-tracking
# sheet, then leave the selection where the user's cursor ended up (C10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 -> Monday 2018-11-26, worked 17:20 - 20:30.
# Copy A3's date formatting down to A8 first (so the new date cell picks up
# the same "m/d/yyyy" number format as the rest of the Day column), then
# overwrite the value.
$ws.Range("A3").Copy($ws.Range("A8"))
$ws.Range("A8").Value = 43430
$ws.Range("B8").Value = 0.72222222222222221
$ws.Range("C8").Value = 0.85416666666666663

# Row 9 -> Tuesday 2018-11-27, worked 10:00 - 21:20.
$ws.Range("A3").Copy($ws.Range("A9"))
$ws.Range("A9").Value = 43431
$ws.Range("B9").Value = 0.41666666666666669
$ws.Range("C9").Value = 0.88888888888888884

# The D column already carries the shared "=Cn-Bn" formula down through row
# 14, and E3's "=SUM(D3:D100)" total will pick up the new rows automatically
# on recalculation.

# Move the selection to where it ended up after entering the new rows.
$ws.Range("C10").Select() | Out-Null
